$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update the Date value ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-09-12T13:34:32+00:00"

# --- Concepts sheet: add "Title Case" display variants for a few rows ---
$concepts = $wb.Worksheets.Item("Concepts")
$concepts.Range("C3").Value = "Missing - Restricted Access"
$concepts.Range("C4").Value = "Missing - Not Provided"
$concepts.Range("C5").Value = "Missing - Not Collected"
$concepts.Range("C7").Value = "Prefer not to Answer"
